$wb = $excel.ActiveWorkbook

# Add the new "Credentials" sheet after the last existing sheet ("Validation part")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Credentials"

# Fill in data in the exact order it was typed (bottom row to top row, right to left
# within a row) so the shared-string table is interned in the original order.
$ws.Range("F8").Value = 0
$ws.Range("E8").Value = "/cluster/storage/no-backup/ccn/CcnStorage1/CCNCDR44/archive/"
$ws.Range("D8").Value = "CCNtasuser@123"
$ws.Range("C8").Value = "tasuser"
$ws.Range("B8").Value = "10.95.213.132"
$ws.Range("A8").Value = "CCN1"

$ws.Range("F7").Value = 90
$ws.Range("E7").Value = "/cluster/storage/no-backup/ccn/CcnStorage0/CCNCDR44/archive/"
$ws.Range("D7").Value = "CCNtasuser@123"
$ws.Range("C7").Value = "tasuser"
$ws.Range("B7").Value = "10.95.213.132"
$ws.Range("A7").Value = "CCN0"

$ws.Range("F6").Value = 60
$ws.Range("E6").Value = "/var/opt/air/datarecords/backup_CDR/"
$ws.Range("D6").Value = "Ericssondu@123"
$ws.Range("C6").Value = "tasuser"
$ws.Range("B6").Value = "10.95.214.166"
$ws.Range("A6").Value = "AIR"

$ws.Range("F5").Value = 0
$ws.Range("E5").Value = "/home/tasuser"
$ws.Range("D5").Value = "Ericssondu@123"
$ws.Range("C5").Value = "tasuser"
$ws.Range("B5").Value = "10.95.214.22"
$ws.Range("A5").Value = "OCC2"

$ws.Range("F4").Value = 150
$ws.Range("E4").Value = "/home/tasuser"
$ws.Range("D4").Value = "Ericssondu@123"
$ws.Range("C4").Value = "tasuser"
$ws.Range("B4").Value = "10.95.214.21"
$ws.Range("A4").Value = "OCC1"

$ws.Range("F3").Value = 150
$ws.Range("E3").Value = "/var/opt/fds/CDR/archive/"
$ws.Range("D3").Value = "Ericssondu@123"
$ws.Range("C3").Value = "tasuser"
$ws.Range("B3").Value = "10.95.214.6"
$ws.Range("A3").Value = "SDP"

$ws.Range("F2").Value = 30
$ws.Range("E2").Value = "/data/fdp/logs/defaultCircle"
$ws.Range("D2").Value = "VenuReddyGaddam"
$ws.Range("C2").Value = "VenuReddyGaddam"
$ws.Range("B2").Value = "10.95.214.72"
$ws.Range("A2").Value = "CIS"

$ws.Range("F1").Value = "Wait_Time"
$ws.Range("E1").Value = "Path"
$ws.Range("D1").Value = "Password"
$ws.Range("C1").Value = "User_Name"
$ws.Range("B1").Value = "IP_HostName"
$ws.Range("A1").Value = "Unix_System"

# Header row shading (new themed fill)
$headerRange = $ws.Range("A1:F1")
$headerRange.Interior.ThemeColor = 5
$headerRange.Interior.TintAndShade = 0.6

# Column widths (best-fit look matching the authored sheet)
$ws.Columns("A").ColumnWidth = 12.43
$ws.Columns("B").ColumnWidth = 13.14
$ws.Columns("C:D").ColumnWidth = 18.86
$ws.Columns("E").ColumnWidth = 60.86
$ws.Columns("F").ColumnWidth = 10.57

$ws.Range("G6").Select()

# Return focus to the originally active sheet/cell
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Range("F2").Select()
